# Apply cell value updates per the crypto price refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure every touched cell is forced to Text format so that numeric-looking
# strings (e.g. "1.00", "0.190", "31.50") keep their exact original formatting
# instead of being auto-coerced into numbers by Excel.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '69.135.96'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -1.50%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.516.59'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -1.95%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '571.31'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -1.44%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '181.69'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -4.63%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.516.76'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -1.84%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.612'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -3.35%  '
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -0.06%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.190'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +5.81%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -3.85%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '53.61'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -5.40%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -0.38%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.44'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -3.67%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.097.59'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -1.74%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '19.23'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -4.24%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.520.56'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '69.085.16'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -1.67%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.43'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -1.08%  '
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -1.30%  '
$ws.Range('B21').NumberFormat = '@'
$ws.Range('B21').Value = 'Polygon'
$ws.Range('C21').NumberFormat = '@'
$ws.Range('C21').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.03'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -1.60%  '
$ws.Range('B22').NumberFormat = '@'
$ws.Range('B22').Value = 'BitcoinCash'
$ws.Range('C22').NumberFormat = '@'
$ws.Range('C22').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '530.43'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +11.91%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '19.79'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +2.39%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.95'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -3.25%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.36'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.35%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '93.92'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +5.58%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.96'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -1.42%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -5.79%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.07'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -3.74%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '31.50'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -2.22%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.26'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -5.82%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '12.54'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +3.08%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '64.35'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -2.55%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.113'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -6.19%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '569.51'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -3.15%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.09'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +6.32%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '37.95'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -3.85%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.00'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -0.07%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -0.88%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -6.45%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.35'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -6.13%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -7.11%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.04'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -5.34%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.53'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +5.16%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.96'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -5.83%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0441'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -1.20%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.161.67'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -2.40%  '
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -4.03%  '
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -3.20%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.996'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -0.20%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '136.57'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -1.21%  '
